# New integral-timing run ("new plots after new integral slurm"): refresh
# the raw For Loop / Numpy / Numexpr timing measurements in B2:D5.
# Columns E:G hold formulas (=$B$2/B2, etc.) that recompute automatically
# off these raw values, and the six embedded scatter charts are bound to
# Sheet1!$B$2:$D$5 so their displayed series track the same cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 358.37575099999998
$ws.Range("C2").Value = 15.571873999999999
$ws.Range("D2").Value = 13.190033

$ws.Range("B3").Value = 348.74953399999998
$ws.Range("C3").Value = 15.549371000000001
$ws.Range("D3").Value = 8.6977700000000002

$ws.Range("B4").Value = 360.318715
$ws.Range("C4").Value = 15.671035
$ws.Range("D4").Value = 6.5458550000000004

$ws.Range("B5").Value = 369.46292
$ws.Range("C5").Value = 15.607931000000001
$ws.Range("D5").Value = 5.1077979999999998

$excel.CalculateFullRebuild()
$wb.Save()
